$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (Giorno, TISG, fcs, buy, MYDIR, need_to_buy)
$data = @(
    @(2,  45985, 11855.5201059538,  10747.4722735802,  17883.4,   6982.5121524095,  -6.39231558376256),
    @(3,  45986, 11246.1959843256,  10227.425817011,   9003.4,    6545.67286428149,  323.737445053855),
    @(4,  45987, 11274.4289271456,  10521.055666338,   9003.4,    6596.42862980057,  338.08684567244),
    @(5,  45988, 11448.556477391,   10857.8094195127,  9003.4,    6758.02533398951,  358.851448062593),
    @(6,  45989, 11505.8318320998,  10268.7571964524,  9003.4,    6804.48193400588,  336.243297102431),
    @(7,  45990, 4441.27964694584,  7133.52972940191,  9003.4,    6031.6631998712,   173.408038719713),
    @(8,  45991, 4285.59691763933,  6968.92917610157,  9003.4,    6023.17493479135,  166.196004620538),
    @(9,  45992, 11775.3234136769,  10665.406686807,   8664.26,   7488.48046616165,  395.401131373692),
    @(10, 45993, 11775.3234136769,  10728.8312027246,  8664.26,   7488.48046616165,  398.043819536926),
    @(11, 45994, 11775.3234136769,  10738.5276501508,  8664.26,   7488.48046616165,  398.447838179684),
    @(12, 45995, 11775.3234136769,  10715.4907220434,  8664.26,   7488.48046616165,  397.487966175211),
    @(13, 45996, 11775.3234136769,  9919.85883698321,  8664.26,   7488.48046616165,  364.336637631036),
    @(14, 45997, 4605.75571869631,  6588.09592649185,  8664.26,   7098.30328362905,  209.255800421704),
    @(15, 45998, 4450.60726673227,  6073.64190747723,  8664.26,   7083.23595700977,  187.192411020292)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
